$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add end time for row 3 (C3), matching the existing date/time formatting
$ws.Range("C3").Value = 42971.958333333336
$ws.Range("C3").NumberFormat = "m/d/yy h:mm"

# Add start time for row 4 (B4), matching the existing date/time formatting
$ws.Range("B4").Value = 42973.395833333336
$ws.Range("B4").NumberFormat = "m/d/yy h:mm"

# Move the active selection to B4
$ws.Range("B4").Select()
